$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Importance values (column B) for rows 58-113
$ws.Range("B58").Value = 59.408000000000001
$ws.Range("B59").Value = 143.13399999999999
$ws.Range("B60").Value = 123.791
$ws.Range("B61").Value = 91.686999999999998
$ws.Range("B62").Value = 133.501
$ws.Range("B63").Value = 146.68199999999999
$ws.Range("B64").Value = 78.361999999999995
$ws.Range("B65").Value = 81.802000000000007
$ws.Range("B66").Value = 78.174000000000007
$ws.Range("B67").Value = 79.5
$ws.Range("B68").Value = 141.53800000000001
$ws.Range("B69").Value = 85.094999999999999
$ws.Range("B70").Value = 99.381
$ws.Range("B71").Value = 74.054000000000002
$ws.Range("B72").Value = 83.045000000000002
$ws.Range("B73").Value = 321.85599999999999
$ws.Range("B74").Value = 89.587999999999994
$ws.Range("B75").Value = 101.408
$ws.Range("B76").Value = 82.043000000000006
$ws.Range("B77").Value = 93.543999999999997
$ws.Range("B78").Value = 166.53
$ws.Range("B79").Value = 82.834999999999994
$ws.Range("B80").Value = 95.105000000000004
$ws.Range("B81").Value = 91.308000000000007
$ws.Range("B82").Value = 104.496
$ws.Range("B83").Value = 168.26400000000001
$ws.Range("B84").Value = 79.664000000000001
$ws.Range("B85").Value = 83.14
$ws.Range("B86").Value = 81.328999999999994
$ws.Range("B87").Value = 73.867999999999995
$ws.Range("B88").Value = 128.28399999999999
$ws.Range("B89").Value = 136.453
$ws.Range("B90").Value = 95.382999999999996
$ws.Range("B91").Value = 90.950999999999993
$ws.Range("B92").Value = 88.119
$ws.Range("B93").Value = 89.825999999999993
$ws.Range("B94").Value = 72.498000000000005
$ws.Range("B95").Value = 65.227999999999994
$ws.Range("B96").Value = 65.912999999999997
$ws.Range("B97").Value = 77.915999999999997
$ws.Range("B98").Value = 62.484999999999999
$ws.Range("B99").Value = 71.840999999999994
$ws.Range("B100").Value = 60.262999999999998
$ws.Range("B101").Value = 61.933
$ws.Range("B102").Value = 106.29900000000001
$ws.Range("B103").Value = 70.507999999999996
$ws.Range("B104").Value = 112.47
$ws.Range("B105").Value = 75.475999999999999
$ws.Range("B106").Value = 78.953999999999994
$ws.Range("B107").Value = 61.826999999999998
$ws.Range("B108").Value = 71.031999999999996
$ws.Range("B109").Value = 62.521000000000001
$ws.Range("B110").Value = 62.558999999999997
$ws.Range("B111").Value = 63.462000000000003
$ws.Range("B112").Value = 66.397000000000006
$ws.Range("B113").Value = 60.197000000000003

# Update sheet view: scroll position, zoom, and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 67
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H88").Select()
$excel.ActiveWindow.Zoom = 70
